$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title / label text (EIA monthly refresh: October 2016 -> November 2016) ---
$ws.Cells.Replace("by Sector, 2006-October 2016 (Thousand Tons)", "by Sector, 2006-November 2016 (Thousand Tons)")
$ws.Cells.Replace("Rolling 12 Months Ending in October", "Rolling 12 Months Ending in November")

# --- Insert the new "November" monthly row under the Year-2016 block (old row 53 -> new row 54, etc.) ---
$ws.Rows(53).Insert()

# Copy formatting from the existing "October" month row (row 52) down onto the newly inserted row,
# so the new row picks up the same styles (s=8 / s=9) used by every other month row.
$ws.Range("A52:F52").Copy()
$ws.Range("A53").PasteSpecial(-4122)

$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 1036
$ws.Range("C53").Value = 84
$ws.Range("D53").Value = 120
$ws.Range("E53").Value = 48
$ws.Range("F53").Value = 784

# --- Refresh the "Year to Date" annual rows (now shifted to rows 55-57) ---
$ws.Range("B55").Value = 16534
$ws.Range("C55").Value = 891
$ws.Range("D55").Value = 1679
$ws.Range("E55").Value = 787
$ws.Range("F55").Value = 13177

$ws.Range("B56").Value = 15269
$ws.Range("C56").Value = 937
$ws.Range("D56").Value = 1829
$ws.Range("E56").Value = 577
$ws.Range("F56").Value = 11927

$ws.Range("B57").Value = 13016
$ws.Range("C57").Value = 919
$ws.Range("D57").Value = 1594
$ws.Range("E57").Value = 488
$ws.Range("F57").Value = 10016

# --- Refresh the "Rolling 12 Months" annual rows (now shifted to rows 59-60) ---
$ws.Range("B59").Value = 16842
$ws.Range("C59").Value = 1023
$ws.Range("D59").Value = 1971
$ws.Range("E59").Value = 650
$ws.Range("F59").Value = 13197

$ws.Range("B60").Value = 14379
$ws.Range("C60").Value = 1014
$ws.Range("D60").Value = 1745
$ws.Range("E60").Value = 546
$ws.Range("F60").Value = 11074
